$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.487.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.035.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.47"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.38"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.033.97"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -14.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.483"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.27"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.533.94"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.599.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.035.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.109"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.02"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.678"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.47"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.48"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.42"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.04"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.01"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.97"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.42"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.64"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.16"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.62"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0405"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "438.53"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0807"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.990.50"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.266"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.52"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.26%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.93"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₃0508"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.06"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.45%  "
